$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 260 — pushes the existing rows 260..339 down to 261..340.
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with the new weekly record.
$ws.Range("A260").Value = 7
$ws.Range("B260").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C260").Value = "Ñuble"
$ws.Range("D260").Value = 44627
$ws.Range("E260").Value = 16
$ws.Range("F260").Value = "Fruta"
$ws.Range("G260").Value = 100102
$ws.Range("H260").Value = "Cítricos"
$ws.Range("I260").Value = 100102005
$ws.Range("J260").Value = "Naranja"
$ws.Range("K260").Value = "Valencia"
$ws.Range("L260").Value = "Primera"
$ws.Range("M260").Value = 120
$ws.Range("N260").Value = 10000
$ws.Range("O260").Value = 11000
$ws.Range("P260").Value = 10500
$ws.Range("Q260").Value = "$/bandeja 15 kilos granel"
$ws.Range("R260").Value = "Región de O'Higgins"
$ws.Range("S260").Value = 700
$ws.Range("T260").Value = 15

Write-Output "row 260 inserted and populated"
